# Edit script implementing the XML diff:
# 1. Insert "the " before the first "date of birth" (paragraph 1 - "requesting the correction of date of birth")
# 2. Change "Correction Of Date Of Birth" -> "Correction of Date of Birth", split into
#    3 runs and insert a collapsed "_GoBack" bookmark between run 2 and run 3.
# 3. Insert a "Number" run after the "ippis" run in "IPPIS No: {ippis}" (table cell).
# 4. Remove the stand-alone "_GoBack" bookmark further down the document (its own empty paragraph).
#    (Bookmark ids auto-renumber by document position, which reproduces the id 0/1 swap for
#    the newly created _GoBack bookmark and the existing _Hlk111541043 bookmark.)
# 5. Insert a "Number" run after the "{ippis" run (second IPPIS placeholder).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. " requesting the correction of " + "date of birth" -> insert "the " run
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Start = 0
$rng.Find.Execute("date of birth", $true, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$insPos = $rng.Start
$insRng = $d.Range($insPos, $insPos)
$insRng.InsertBefore("the ")
# Force the new runs to stay split from their neighbours (identical formatting
# would otherwise cause them to be coalesced back together) by toggling a
# character property on and back off at each boundary.
$splitRng = $d.Range($insPos, $insPos + 4)
$splitRng.Font.Bold = 1
$splitRng.Font.Bold = 0
$splitRng2 = $d.Range($insPos + 4, $insPos + 17)
$splitRng2.Font.Bold = 1
$splitRng2.Font.Bold = 0

Write-Output "Step 1 done"

# ---------------------------------------------------------------------------
# 2. "Correction Of Date Of Birth" -> "Correction of Date of Birth" split into
#    three runs ("Correction o" / "f Date o" / "f Birth") with a collapsed
#    "_GoBack" bookmark inserted between run 2 and run 3.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Start = 0
$rng2.Find.Execute("Correction Of Date Of Birth", $true, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$hStart = $rng2.Start
$hRng = $d.Range($hStart, $hStart + 28)
$hRng.Text = "Correction of Date of Birth"

# Force run splits at offsets 12 and 20 (relative to $hStart) by toggling Bold.
$split1 = $d.Range($hStart, $hStart + 12)
$split1.Font.Bold = 0
$split1.Font.Bold = 1
$split2 = $d.Range($hStart + 12, $hStart + 20)
$split2.Font.Bold = 0
$split2.Font.Bold = 1
$split3 = $d.Range($hStart + 20, $hStart + 27)
$split3.Font.Bold = 0
$split3.Font.Bold = 1

# Insert the collapsed "_GoBack" bookmark right at the boundary between
# run 2 ("f Date o") and run 3 ("f Birth").
$bmRng = $d.Range($hStart + 20, $hStart + 20)
$d.Bookmarks.Add("_GoBack", $bmRng)

Write-Output "Step 2 done"

# ---------------------------------------------------------------------------
# 3. "IPPIS No: {ippis}" -> "IPPIS No: {ippisNumber}" (insert "Number" run
#    right before the closing "}").
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Start = 0
$rng3.Find.Execute("IPPIS No: {ippis}", $true, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$s3 = $rng3.Start
$insPos3 = $s3 + 16
$insRng3 = $d.Range($insPos3, $insPos3)
$insRng3.InsertBefore("Number")
$splitRng3 = $d.Range($insPos3, $insPos3 + 6)
$splitRng3.Font.Bold = 1
$splitRng3.Font.Bold = 0

Write-Output "Step 3 done"

# ---------------------------------------------------------------------------
# 4. "IPPIS N0: {ippis}" -> "IPPIS N0: {ippisNumber}" (insert "Number" run
#    right before the closing "}").
# ---------------------------------------------------------------------------
$rng4 = $d.Content
$rng4.Start = 0
$rng4.Find.Execute("IPPIS N0: {ippis}", $true, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$s4 = $rng4.Start
$insPos4 = $s4 + 16
$insRng4 = $d.Range($insPos4, $insPos4)
$insRng4.InsertBefore("Number")
$splitRng4 = $d.Range($insPos4, $insPos4 + 6)
$splitRng4.Font.Bold = 0
$splitRng4.Font.Bold = 1

Write-Output "Step 4 done"
